# The post "「レットイットゴー／let it go」" (previously row 857) was removed.
# All subsequent rows shift up by one, and the used range shrinks from
# A1:C867 to A1:C866. Deleting the entire row in Excel reproduces exactly
# this shift (cells below move up, dimension auto-updates).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(857).Delete()
